# Logged Week 16 and performed season sim from Week 17
# Update the "R" (Road) row (row 3) target depth totals on both the
# OFF and DEF sheets to reflect the newly logged week.

$wb = $excel.ActiveWorkbook

# OFF sheet - row 3 ("R")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 187
$wsOff.Range("C3").Value = 130
$wsOff.Range("D3").Value = 50
$wsOff.Range("E3").Value = 29
$wsOff.Range("F3").Value = 3

# DEF sheet - row 3 ("R")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 166
$wsDef.Range("C3").Value = 132
$wsDef.Range("D3").Value = 41
$wsDef.Range("E3").Value = 19
